# expressionRuntimeException-expected-generation.docx was regenerated against a
# newer build (Acceleo Query Evaluator / M2DocEvaluator / JUnit4 runner), which
# shifted several source line numbers recorded in the embedded Java stack
# trace ("Moving from 3.2.0 to 3.2.1"). Rebuild the old/new stack-trace text
# (tab-indented, LF-separated, exactly as stored in the run) and swap it in a
# single Find & Replace so every affected "xxx.java:NNN)" reference - plus the
# JUnit ParentRunner call chain that was restructured - is updated at once.

$oldLines = @(
    "`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)", "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)", "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:186)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:342)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:262)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)",
    "`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1296)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:147)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)",
    "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)", "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)", "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)", "`tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)", "`tat java.lang.reflect.Method.invoke(Method.java:498)", "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)", "`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)",
    "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)", "`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)", "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)", "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)", "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)", "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)", "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)", "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)", "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)", "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)", "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", "`tat org.junit.runners.Suite.runChild(Suite.java:128)", "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)", "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)", "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)", "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)", "`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)", "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)", "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)",
    "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", "`tat org.junit.runners.Suite.runChild(Suite.java:128)", "`tat org.junit.runners.Suite.runChild(Suite.java:27)", "`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)",
    "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)", "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)", "`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)", "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)",
    "`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)", "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)", "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)", "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)"
)

$newLines = @(
    "`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)", "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)", "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseSequenceInExtensionLiteral(AstEvaluator.java:345)",
    "`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:262)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)",
    "`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1415)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseRepetition(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:147)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)", "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)", "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)", "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)",
    "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)", "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)", "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)", "`tat sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)",
    "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)", "`tat java.lang.reflect.Method.invoke(Method.java:498)", "`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)", "`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)",
    "`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)", "`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)", "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)", "`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)",
    "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)", "`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)", "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)", "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)", "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)", "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)", "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)",
    "`tat org.junit.runners.Suite.runChild(Suite.java:128)", "`tat org.junit.runners.Suite.runChild(Suite.java:27)", "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)", "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)",
    "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)", "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)", "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)", "`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)",
    "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)", "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", "`tat org.junit.runners.Suite.runChild(Suite.java:128)", "`tat org.junit.runners.Suite.runChild(Suite.java:27)",
    "`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)", "`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)", "`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)", "`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)",
    "`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)", "`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)", "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)", "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)",
    "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)"
)

$old = [string]::Join("`n", $oldLines)
$new = [string]::Join("`n", $newLines)

$d = $word.ActiveDocument
$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the expected stack trace block to replace."
}

Write-Output ("Replaced stack trace block: " + $found)
